$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B6").Value = "SingleUseId3"
$ws.Range("C6").Value = "Default"
$ws.Range("D6").Value = "Center"
$ws.Range("E6").Value = "LTR"
$ws.Range("F6").Value = "Resource 2"
